# Append a new "2025-04-02" row (row 32) to every price sheet, carrying
# forward the most recent (row 31) price for that sheet — matching the
# commit "Updated Argent prices in Excel".
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Grab the latest known price as literal text (Value2 avoids the
    # COM shim's broken Value getter) before we touch formatting.
    $priceValue = $ws.Range("B31").Value2

    # Force the new cells to Text so "2025-04-02" and values like
    # "5,489" are stored as literal strings instead of being
    # auto-coerced into a date serial / thousands-formatted number.
    $ws.Range("A32:B32").NumberFormat = "@"
    $ws.Range("A32").Value = "2025-04-02"
    $ws.Range("B32").Value = $priceValue

    # Drop the temporary Text format off the cells again so the new
    # row matches the plain (unstyled) look of the existing rows.
    $ws.Range("A32:B32").ClearFormats()
}
